$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "MasterProtected"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("MasterProtected")

# Row 2: Source IP now holds two comma separated CIDRs (wrapped across two
# lines) and Protocol switches from TCP to udp.
$ws1.Range("C2").Value = "10.108.0.0/15," + [char]10 + "10.110.64.0/22"
$ws1.Range("C2").WrapText = $true
$ws1.Range("F2").Value = "udp"
$ws1.Rows.Item(2).RowHeight = 28.8

# Row 4 was a stray blank row in the template - remove it entirely so the
# row node disappears from the sheet (contents + formatting).
$ws1.Range("A4:J4").Clear()

$ws1.Range("D5").Select()

# ---------------------------------------------------------------------
# Sheet "FirewallRulesToValidate"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("FirewallRulesToValidate")

# Row 2: Source IP / Destination IP now hold comma separated values.
$ws2.Range("B2").Value = "200.0.32.0/22," + [char]10 + "10.110.64.1"
$ws2.Range("B2").WrapText = $true
$ws2.Range("D2").Value = "20.0.0.1,30.0.0.1"
$ws2.Rows.Item(2).RowHeight = 28.8

# New rows for the additional rules being validated.
$ws2.Range("A5").Value = "Azure Non prod"
$ws2.Range("B5").Value = "10.110.64.0/22"
$ws2.Range("C5").Value = "IOD- Non-Prod"
$ws2.Range("D5").Value = "30.0.0.1/32"
$ws2.Range("E5").Value = "TCP"
$ws2.Range("F5").Value = 443
$ws2.Range("G5").Value = "HTTPS"

$ws2.Range("A6").Value = "Azure Non prod"
$ws2.Range("B6").Value = "10.110.64.0/22"
$ws2.Range("C6").Value = "IOD- Non-Prod"
$ws2.Range("D6").Value = "30.0.0.1/32"
$ws2.Range("E6").Value = "ICMP"
$ws2.Range("G6").Value = "icmp"

$ws2.Range("E7").Select()
